$p = $ppt.ActivePresentation

# Slide 32 -> Title "HTML - Tags Sem Fechamento"
$s = $p.Slides.Item(32)
$titleShape = $s.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

$origLen = $tr.Length

# Append a new run containing " *" right after the existing title text
# (i.e. after the "Fechamento" run), mirroring the target OOXML which adds
# a brand new <a:r> run before the paragraph's endParaRPr.
$null = $tr.InsertAfter(" *")

# Scope formatting to just the newly inserted text so the preceding runs
# ("HTML - Tags Sem " / "Fechamento") are left untouched.
$newRun = $tr.Characters($origLen + 1, 2)
$newRun.Font.Bold = $true
$newRun.Font.Color.RGB = 12611584   # 0x0070C0 (stored as BGR by the Font.Color API)
